$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-25)
# from serial date 45236 to 45237 (one day later).
for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45236) {
        $cell.Value2 = 45237
    }
}
